# Update error rate values on the ErrorRates worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ErrorRates")

# Row 2 (User 1)
$ws.Range("D2").Value = 0.05

# Row 5 (User 4)
$ws.Range("B5").Value = 0.2
$ws.Range("D5").Value = 0.1

# Row 13 (User 12)
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0

# Row 14 (User 13)
$ws.Range("B14").Value = 0.1
$ws.Range("D14").Value = 0.05

# Row 20 (User 19)
$ws.Range("D20").Value = 0

# Row 24 (User 23)
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0.05

# Row 29 (User 28)
$ws.Range("B29").Value = 0.2
$ws.Range("D29").Value = 0.1
